# "changed fragment duration test. updated results"
#
# The author reran the low-latency fragment-duration experiment: several
# measured "latency" values (column C) were corrected, a new "Buffer
# Length" metric (column H) was added with a header and a couple of data
# points, and a few rows that previously had no latency reading now have
# one. Finally the selection cursor was left on T9 when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header: H1 = "Buffer Length"
$ws.Range("H1").Value = "Buffer Length"

# Corrected / newly-recorded latency readings (column C)
$ws.Range("C10").Value = 7.5    # was 6.5
$ws.Range("C11").Value = 9      # was 10
$ws.Range("C12").Value = 12     # was 11
$ws.Range("C13").Value = 15     # previously blank
$ws.Range("C14").Value = 20     # previously blank
$ws.Range("C15").Value = 25     # previously blank

# New Buffer Length measurements (column H)
$ws.Range("H14").Value = 0.2
$ws.Range("H15").Value = 0.1

# Leave the selection where the author left it when saving
$ws.Range("T9").Select()
